# Sincronizando o repositório local baixado do Luciano com o repositório que eu criei
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 8422
$ws.Range("B2").Value = "Caroline Castro"
$ws.Range("C2").Value = "Operações"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45090
$ws.Range("G2").Value = 11650.51

# Row 3
$ws.Range("A3").Value = 10445
$ws.Range("B3").Value = "Diego Rezende"
$ws.Range("C3").Value = "Operações"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45093
$ws.Range("G3").Value = 2804.65

# Row 4
$ws.Range("A4").Value = 53933
$ws.Range("B4").Value = "Srta. Laís Nogueira"
$ws.Range("C4").Value = "Marketing"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45085
$ws.Range("G4").Value = 4989.95

# Row 5
$ws.Range("A5").Value = 13402
$ws.Range("B5").Value = "João Costa"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 45080
$ws.Range("G5").Value = 12360.44

# Row 6
$ws.Range("A6").Value = 79511
$ws.Range("B6").Value = "Antônio Fernandes"
$ws.Range("C6").Value = "Jurídico"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45092
$ws.Range("G6").Value = 5412.64

# Row 7
$ws.Range("A7").Value = 85504
$ws.Range("B7").Value = "Samuel Almeida"
$ws.Range("C7").Value = "Jurídico"
$ws.Range("D7").Value = "Viagem de negócios"
$ws.Range("F7").Value = 45094
$ws.Range("G7").Value = 6319.1

# Row 8
$ws.Range("A8").Value = 2214
$ws.Range("B8").Value = "Kamilly Mendes"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45088
$ws.Range("G8").Value = 9609.48

# Row 9
$ws.Range("A9").Value = 23130
$ws.Range("B9").Value = "Luiz Miguel da Cunha"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45101
$ws.Range("G9").Value = 7114.54

# Row 10
$ws.Range("A10").Value = 60738
$ws.Range("B10").Value = "Vitor Aragão"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45096
$ws.Range("G10").Value = 7532.55

# Row 11
$ws.Range("A11").Value = 81395
$ws.Range("B11").Value = "Marcelo Araújo"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45080
$ws.Range("G11").Value = 3726.96
